# === Edit script: add "Criminal Offenses - Rape" rows (replicating prior commit pattern) ===
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Populate the new data rows (578-649) ---
#     Offense = "Criminal Offenses - Rape", Date in {sum2014, sum2015}
#     across Reporting Location x Sector of Institution combinations.

$ws.Cells.Item(578,1).Value = 'Public, 4-year or above'
$ws.Cells.Item(578,2).Value = 'On Campus (excluding Residence Halls)'
$ws.Cells.Item(578,3).Value = 'Criminal Offenses - Rape'
$ws.Cells.Item(578,4).Value = 'sum2014'
$ws.Cells.Item(578,5).Value = 426
$ws.Cells.Item(579,1).Value = 'Private nonprofit, 4-year or above'
$ws.Cells.Item(579,2).Value = 'On Campus (excluding Residence Halls)'
$ws.Cells.Item(579,3).Value = 'Criminal Offenses - Rape'
$ws.Cells.Item(579,4).Value = 'sum2014'
$ws.Cells.Item(579,5).Value = 315
$ws.Cells.Item(580,1).Value = 'Private for-profit, 4-year or above'
$ws.Cells.Item(580,2).Value = 'On Campus (excluding Residence Halls)'
$ws.Cells.Item(580,3).Value = 'Criminal Offenses - Rape'
$ws.Cells.Item(580,4).Value = 'sum2014'
$ws.Cells.Item(580,5).Value = 5
$ws.Cells.Item(581,1).Value = 'Public, 2-year'
$ws.Cells.Item(581,2).Value = 'On Campus (excluding Residence Halls)'
$ws.Cells.Item(581,3).Value = 'Criminal Offenses - Rape'
$ws.Cells.Item(581,4).Value = 'sum2014'
$ws.Cells.Item(581,5).Value = 57
$ws.Cells.Item(582,1).Value = 'Private nonprofit, 2-year'
$ws.Cells.Item(582,2).Value = 'On Campus (excluding Residence Halls)'
$ws.Cells.Item(582,3).Value = 'Criminal Offenses - Rape'
$ws.Cells.Item(582,4).Value = 'sum2014'
$ws.Cells.Item(582,5).Value = 2
$ws.Cells.Item(583,1).Value = 'Private for-profit, 2-year'
$ws.Cells.Item(583,2).Value = 'On Campus (excluding Residence Halls)'
$ws.Cells.Item(583,3).Value = 'Criminal Offenses - Rape'
$ws.Cells.Item(583,4).Value = 'sum2014'
$ws.Cells.Item(583,5).Value = 0
$ws.Cells.Item(584,1).Value = 'Public, less-than 2-year'
$ws.Cells.Item(584,2).Value = 'On Campus (excluding Residence Halls)'
$ws.Cells.Item(584,3).Value = 'Criminal Offenses - Rape'
$ws.Cells.Item(584,4).Value = 'sum2014'
$ws.Cells.Item(584,5).Value = 0
$ws.Cells.Item(585,1).Value = 'Private nonprofit, less-than 2-year'
$ws.Cells.Item(585,2).Value = 'On Campus (excluding Residence Halls)'
$ws.Cells.Item(585,3).Value = 'Criminal Offenses - Rape'
$ws.Cells.Item(585,4).Value = 'sum2014'
$ws.Cells.Item(585,5).Value = 0
$ws.Cells.Item(586,1).Value = 'Private for-profit, less-than 2-year'
$ws.Cells.Item(586,2).Value = 'On Campus (excluding Residence Halls)'
$ws.Cells.Item(586,3).Value = 'Criminal Offenses - Rape'
$ws.Cells.Item(586,4).Value = 'sum2014'
$ws.Cells.Item(586,5).Value = 1
$ws.Cells.Item(587,1).Value = 'Public, 4-year or above'
$ws.Cells.Item(587,2).Value = 'On Campus (excluding Residence Halls)'
$ws.Cells.Item(587,3).Value = 'Criminal Offenses - Rape'
$ws.Cells.Item(587,4).Value = 'sum2015'
$ws.Cells.Item(587,5).Value = 481
$ws.Cells.Item(588,1).Value = 'Private nonprofit, 4-year or above'
$ws.Cells.Item(588,2).Value = 'On Campus (excluding Residence Halls)'
$ws.Cells.Item(588,3).Value = 'Criminal Offenses - Rape'
$ws.Cells.Item(588,4).Value = 'sum2015'
$ws.Cells.Item(588,5).Value = 352
$ws.Cells.Item(589,1).Value = 'Private for-profit, 4-year or above'
$ws.Cells.Item(589,2).Value = 'On Campus (excluding Residence Halls)'
$ws.Cells.Item(589,3).Value = 'Criminal Offenses - Rape'
$ws.Cells.Item(589,4).Value = 'sum2015'
$ws.Cells.Item(589,5).Value = 3
$ws.Cells.Item(590,1).Value = 'Public, 2-year'
$ws.Cells.Item(590,2).Value = 'On Campus (excluding Residence Halls)'
$ws.Cells.Item(590,3).Value = 'Criminal Offenses - Rape'
$ws.Cells.Item(590,4).Value = 'sum2015'
$ws.Cells.Item(590,5).Value = 69
$ws.Cells.Item(591,1).Value = 'Private nonprofit, 2-year'
$ws.Cells.Item(591,2).Value = 'On Campus (excluding Residence Halls)'
$ws.Cells.Item(591,3).Value = 'Criminal Offenses - Rape'
$ws.Cells.Item(591,4).Value = 'sum2015'
$ws.Cells.Item(591,5).Value = 2
$ws.Cells.Item(592,1).Value = 'Private for-profit, 2-year'
$ws.Cells.Item(592,2).Value = 'On Campus (excluding Residence Halls)'
$ws.Cells.Item(592,3).Value = 'Criminal Offenses - Rape'
$ws.Cells.Item(592,4).Value = 'sum2015'
$ws.Cells.Item(592,5).Value = 0
$ws.Cells.Item(593,1).Value = 'Public, less-than 2-year'
$ws.Cells.Item(593,2).Value = 'On Campus (excluding Residence Halls)'
$ws.Cells.Item(593,3).Value = 'Criminal Offenses - Rape'
$ws.Cells.Item(593,4).Value = 'sum2015'
$ws.Cells.Item(593,5).Value = 0
$ws.Cells.Item(594,1).Value = 'Private nonprofit, less-than 2-year'
$ws.Cells.Item(594,2).Value = 'On Campus (excluding Residence Halls)'
$ws.Cells.Item(594,3).Value = 'Criminal Offenses - Rape'
$ws.Cells.Item(594,4).Value = 'sum2015'
$ws.Cells.Item(594,5).Value = 0
$ws.Cells.Item(595,1).Value = 'Private for-profit, less-than 2-year'
$ws.Cells.Item(595,2).Value = 'On Campus (excluding Residence Halls)'
$ws.Cells.Item(595,3).Value = 'Criminal Offenses - Rape'
$ws.Cells.Item(595,4).Value = 'sum2015'
$ws.Cells.Item(595,5).Value = 1
$ws.Cells.Item(596,1).Value = 'Public, 4-year or above'
$ws.Cells.Item(596,2).Value = 'On Campus (Residence Halls)'
$ws.Cells.Item(596,3).Value = 'Criminal Offenses - Rape'
$ws.Cells.Item(596,4).Value = 'sum2014'
$ws.Cells.Item(596,5).Value = 1703
$ws.Cells.Item(597,1).Value = 'Private nonprofit, 4-year or above'
$ws.Cells.Item(597,2).Value = 'On Campus (Residence Halls)'
$ws.Cells.Item(597,3).Value = 'Criminal Offenses - Rape'
$ws.Cells.Item(597,4).Value = 'sum2014'
$ws.Cells.Item(597,5).Value = 1858
$ws.Cells.Item(598,1).Value = 'Private for-profit, 4-year or above'
$ws.Cells.Item(598,2).Value = 'On Campus (Residence Halls)'
$ws.Cells.Item(598,3).Value = 'Criminal Offenses - Rape'
$ws.Cells.Item(598,4).Value = 'sum2014'
$ws.Cells.Item(598,5).Value = 23
$ws.Cells.Item(599,1).Value = 'Public, 2-year'
$ws.Cells.Item(599,2).Value = 'On Campus (Residence Halls)'
$ws.Cells.Item(599,3).Value = 'Criminal Offenses - Rape'
$ws.Cells.Item(599,4).Value = 'sum2014'
$ws.Cells.Item(599,5).Value = 77
$ws.Cells.Item(600,1).Value = 'Private nonprofit, 2-year'
$ws.Cells.Item(600,2).Value = 'On Campus (Residence Halls)'
$ws.Cells.Item(600,3).Value = 'Criminal Offenses - Rape'
$ws.Cells.Item(600,4).Value = 'sum2014'
$ws.Cells.Item(600,5).Value = 3
$ws.Cells.Item(601,1).Value = 'Private for-profit, 2-year'
$ws.Cells.Item(601,2).Value = 'On Campus (Residence Halls)'
$ws.Cells.Item(601,3).Value = 'Criminal Offenses - Rape'
$ws.Cells.Item(601,4).Value = 'sum2014'
$ws.Cells.Item(601,5).Value = 1
$ws.Cells.Item(602,1).Value = 'Public, less-than 2-year'
$ws.Cells.Item(602,2).Value = 'On Campus (Residence Halls)'
$ws.Cells.Item(602,3).Value = 'Criminal Offenses - Rape'
$ws.Cells.Item(602,4).Value = 'sum2014'
$ws.Cells.Item(602,5).Value = 0
$ws.Cells.Item(603,1).Value = 'Private nonprofit, less-than 2-year'
$ws.Cells.Item(603,2).Value = 'On Campus (Residence Halls)'
$ws.Cells.Item(603,3).Value = 'Criminal Offenses - Rape'
$ws.Cells.Item(603,4).Value = 'sum2014'
# Row 603: Count intentionally left blank (missing data in source)
$ws.Cells.Item(604,1).Value = 'Private for-profit, less-than 2-year'
$ws.Cells.Item(604,2).Value = 'On Campus (Residence Halls)'
$ws.Cells.Item(604,3).Value = 'Criminal Offenses - Rape'
$ws.Cells.Item(604,4).Value = 'sum2014'
$ws.Cells.Item(604,5).Value = 2
$ws.Cells.Item(605,1).Value = 'Public, 4-year or above'
$ws.Cells.Item(605,2).Value = 'On Campus (Residence Halls)'
$ws.Cells.Item(605,3).Value = 'Criminal Offenses - Rape'
$ws.Cells.Item(605,4).Value = 'sum2015'
$ws.Cells.Item(605,5).Value = 2049
$ws.Cells.Item(606,1).Value = 'Private nonprofit, 4-year or above'
$ws.Cells.Item(606,2).Value = 'On Campus (Residence Halls)'
$ws.Cells.Item(606,3).Value = 'Criminal Offenses - Rape'
$ws.Cells.Item(606,4).Value = 'sum2015'
$ws.Cells.Item(606,5).Value = 2031
$ws.Cells.Item(607,1).Value = 'Private for-profit, 4-year or above'
$ws.Cells.Item(607,2).Value = 'On Campus (Residence Halls)'
$ws.Cells.Item(607,3).Value = 'Criminal Offenses - Rape'
$ws.Cells.Item(607,4).Value = 'sum2015'
$ws.Cells.Item(607,5).Value = 8
$ws.Cells.Item(608,1).Value = 'Public, 2-year'
$ws.Cells.Item(608,2).Value = 'On Campus (Residence Halls)'
$ws.Cells.Item(608,3).Value = 'Criminal Offenses - Rape'
$ws.Cells.Item(608,4).Value = 'sum2015'
$ws.Cells.Item(608,5).Value = 127
$ws.Cells.Item(609,1).Value = 'Private nonprofit, 2-year'
$ws.Cells.Item(609,2).Value = 'On Campus (Residence Halls)'
$ws.Cells.Item(609,3).Value = 'Criminal Offenses - Rape'
$ws.Cells.Item(609,4).Value = 'sum2015'
$ws.Cells.Item(609,5).Value = 1
$ws.Cells.Item(610,1).Value = 'Private for-profit, 2-year'
$ws.Cells.Item(610,2).Value = 'On Campus (Residence Halls)'
$ws.Cells.Item(610,3).Value = 'Criminal Offenses - Rape'
$ws.Cells.Item(610,4).Value = 'sum2015'
$ws.Cells.Item(610,5).Value = 2
$ws.Cells.Item(611,1).Value = 'Public, less-than 2-year'
$ws.Cells.Item(611,2).Value = 'On Campus (Residence Halls)'
$ws.Cells.Item(611,3).Value = 'Criminal Offenses - Rape'
$ws.Cells.Item(611,4).Value = 'sum2015'
$ws.Cells.Item(611,5).Value = 0
$ws.Cells.Item(612,1).Value = 'Private nonprofit, less-than 2-year'
$ws.Cells.Item(612,2).Value = 'On Campus (Residence Halls)'
$ws.Cells.Item(612,3).Value = 'Criminal Offenses - Rape'
$ws.Cells.Item(612,4).Value = 'sum2015'
# Row 612: Count intentionally left blank (missing data in source)
$ws.Cells.Item(613,1).Value = 'Private for-profit, less-than 2-year'
$ws.Cells.Item(613,2).Value = 'On Campus (Residence Halls)'
$ws.Cells.Item(613,3).Value = 'Criminal Offenses - Rape'
$ws.Cells.Item(613,4).Value = 'sum2015'
$ws.Cells.Item(613,5).Value = 1
$ws.Cells.Item(614,1).Value = 'Public, 4-year or above'
$ws.Cells.Item(614,2).Value = 'Non-Campus'
$ws.Cells.Item(614,3).Value = 'Criminal Offenses - Rape'
$ws.Cells.Item(614,4).Value = 'sum2014'
$ws.Cells.Item(614,5).Value = 347
$ws.Cells.Item(615,1).Value = 'Private nonprofit, 4-year or above'
$ws.Cells.Item(615,2).Value = 'Non-Campus'
$ws.Cells.Item(615,3).Value = 'Criminal Offenses - Rape'
$ws.Cells.Item(615,4).Value = 'sum2014'
$ws.Cells.Item(615,5).Value = 130
$ws.Cells.Item(616,1).Value = 'Private for-profit, 4-year or above'
$ws.Cells.Item(616,2).Value = 'Non-Campus'
$ws.Cells.Item(616,3).Value = 'Criminal Offenses - Rape'
$ws.Cells.Item(616,4).Value = 'sum2014'
$ws.Cells.Item(616,5).Value = 12
$ws.Cells.Item(617,1).Value = 'Public, 2-year'
$ws.Cells.Item(617,2).Value = 'Non-Campus'
$ws.Cells.Item(617,3).Value = 'Criminal Offenses - Rape'
$ws.Cells.Item(617,4).Value = 'sum2014'
$ws.Cells.Item(617,5).Value = 12
$ws.Cells.Item(618,1).Value = 'Private nonprofit, 2-year'
$ws.Cells.Item(618,2).Value = 'Non-Campus'
$ws.Cells.Item(618,3).Value = 'Criminal Offenses - Rape'
$ws.Cells.Item(618,4).Value = 'sum2014'
$ws.Cells.Item(618,5).Value = 0
$ws.Cells.Item(619,1).Value = 'Private for-profit, 2-year'
$ws.Cells.Item(619,2).Value = 'Non-Campus'
$ws.Cells.Item(619,3).Value = 'Criminal Offenses - Rape'
$ws.Cells.Item(619,4).Value = 'sum2014'
$ws.Cells.Item(619,5).Value = 1
$ws.Cells.Item(620,1).Value = 'Public, less-than 2-year'
$ws.Cells.Item(620,2).Value = 'Non-Campus'
$ws.Cells.Item(620,3).Value = 'Criminal Offenses - Rape'
$ws.Cells.Item(620,4).Value = 'sum2014'
$ws.Cells.Item(620,5).Value = 0
$ws.Cells.Item(621,1).Value = 'Private nonprofit, less-than 2-year'
$ws.Cells.Item(621,2).Value = 'Non-Campus'
$ws.Cells.Item(621,3).Value = 'Criminal Offenses - Rape'
$ws.Cells.Item(621,4).Value = 'sum2014'
$ws.Cells.Item(621,5).Value = 0
$ws.Cells.Item(622,1).Value = 'Private for-profit, less-than 2-year'
$ws.Cells.Item(622,2).Value = 'Non-Campus'
$ws.Cells.Item(622,3).Value = 'Criminal Offenses - Rape'
$ws.Cells.Item(622,4).Value = 'sum2014'
$ws.Cells.Item(622,5).Value = 2
$ws.Cells.Item(623,1).Value = 'Public, 4-year or above'
$ws.Cells.Item(623,2).Value = 'Non-Campus'
$ws.Cells.Item(623,3).Value = 'Criminal Offenses - Rape'
$ws.Cells.Item(623,4).Value = 'sum2015'
$ws.Cells.Item(623,5).Value = 364
$ws.Cells.Item(624,1).Value = 'Private nonprofit, 4-year or above'
$ws.Cells.Item(624,2).Value = 'Non-Campus'
$ws.Cells.Item(624,3).Value = 'Criminal Offenses - Rape'
$ws.Cells.Item(624,4).Value = 'sum2015'
$ws.Cells.Item(624,5).Value = 150
$ws.Cells.Item(625,1).Value = 'Private for-profit, 4-year or above'
$ws.Cells.Item(625,2).Value = 'Non-Campus'
$ws.Cells.Item(625,3).Value = 'Criminal Offenses - Rape'
$ws.Cells.Item(625,4).Value = 'sum2015'
$ws.Cells.Item(625,5).Value = 5
$ws.Cells.Item(626,1).Value = 'Public, 2-year'
$ws.Cells.Item(626,2).Value = 'Non-Campus'
$ws.Cells.Item(626,3).Value = 'Criminal Offenses - Rape'
$ws.Cells.Item(626,4).Value = 'sum2015'
$ws.Cells.Item(626,5).Value = 16
$ws.Cells.Item(627,1).Value = 'Private nonprofit, 2-year'
$ws.Cells.Item(627,2).Value = 'Non-Campus'
$ws.Cells.Item(627,3).Value = 'Criminal Offenses - Rape'
$ws.Cells.Item(627,4).Value = 'sum2015'
$ws.Cells.Item(627,5).Value = 1
$ws.Cells.Item(628,1).Value = 'Private for-profit, 2-year'
$ws.Cells.Item(628,2).Value = 'Non-Campus'
$ws.Cells.Item(628,3).Value = 'Criminal Offenses - Rape'
$ws.Cells.Item(628,4).Value = 'sum2015'
$ws.Cells.Item(628,5).Value = 1
$ws.Cells.Item(629,1).Value = 'Public, less-than 2-year'
$ws.Cells.Item(629,2).Value = 'Non-Campus'
$ws.Cells.Item(629,3).Value = 'Criminal Offenses - Rape'
$ws.Cells.Item(629,4).Value = 'sum2015'
$ws.Cells.Item(629,5).Value = 0
$ws.Cells.Item(630,1).Value = 'Private nonprofit, less-than 2-year'
$ws.Cells.Item(630,2).Value = 'Non-Campus'
$ws.Cells.Item(630,3).Value = 'Criminal Offenses - Rape'
$ws.Cells.Item(630,4).Value = 'sum2015'
$ws.Cells.Item(630,5).Value = 0
$ws.Cells.Item(631,1).Value = 'Private for-profit, less-than 2-year'
$ws.Cells.Item(631,2).Value = 'Non-Campus'
$ws.Cells.Item(631,3).Value = 'Criminal Offenses - Rape'
$ws.Cells.Item(631,4).Value = 'sum2015'
$ws.Cells.Item(631,5).Value = 2
$ws.Cells.Item(632,1).Value = 'Public, 4-year or above'
$ws.Cells.Item(632,2).Value = 'Public Property'
$ws.Cells.Item(632,3).Value = 'Criminal Offenses - Rape'
$ws.Cells.Item(632,4).Value = 'sum2014'
$ws.Cells.Item(632,5).Value = 74
$ws.Cells.Item(633,1).Value = 'Private nonprofit, 4-year or above'
$ws.Cells.Item(633,2).Value = 'Public Property'
$ws.Cells.Item(633,3).Value = 'Criminal Offenses - Rape'
$ws.Cells.Item(633,4).Value = 'sum2014'
$ws.Cells.Item(633,5).Value = 73
$ws.Cells.Item(634,1).Value = 'Private for-profit, 4-year or above'
$ws.Cells.Item(634,2).Value = 'Public Property'
$ws.Cells.Item(634,3).Value = 'Criminal Offenses - Rape'
$ws.Cells.Item(634,4).Value = 'sum2014'
$ws.Cells.Item(634,5).Value = 11
$ws.Cells.Item(635,1).Value = 'Public, 2-year'
$ws.Cells.Item(635,2).Value = 'Public Property'
$ws.Cells.Item(635,3).Value = 'Criminal Offenses - Rape'
$ws.Cells.Item(635,4).Value = 'sum2014'
$ws.Cells.Item(635,5).Value = 24
$ws.Cells.Item(636,1).Value = 'Private nonprofit, 2-year'
$ws.Cells.Item(636,2).Value = 'Public Property'
$ws.Cells.Item(636,3).Value = 'Criminal Offenses - Rape'
$ws.Cells.Item(636,4).Value = 'sum2014'
$ws.Cells.Item(636,5).Value = 4
$ws.Cells.Item(637,1).Value = 'Private for-profit, 2-year'
$ws.Cells.Item(637,2).Value = 'Public Property'
$ws.Cells.Item(637,3).Value = 'Criminal Offenses - Rape'
$ws.Cells.Item(637,4).Value = 'sum2014'
$ws.Cells.Item(637,5).Value = 16
$ws.Cells.Item(638,1).Value = 'Public, less-than 2-year'
$ws.Cells.Item(638,2).Value = 'Public Property'
$ws.Cells.Item(638,3).Value = 'Criminal Offenses - Rape'
$ws.Cells.Item(638,4).Value = 'sum2014'
$ws.Cells.Item(638,5).Value = 2
$ws.Cells.Item(639,1).Value = 'Private nonprofit, less-than 2-year'
$ws.Cells.Item(639,2).Value = 'Public Property'
$ws.Cells.Item(639,3).Value = 'Criminal Offenses - Rape'
$ws.Cells.Item(639,4).Value = 'sum2014'
$ws.Cells.Item(639,5).Value = 0
$ws.Cells.Item(640,1).Value = 'Private for-profit, less-than 2-year'
$ws.Cells.Item(640,2).Value = 'Public Property'
$ws.Cells.Item(640,3).Value = 'Criminal Offenses - Rape'
$ws.Cells.Item(640,4).Value = 'sum2014'
$ws.Cells.Item(640,5).Value = 5
$ws.Cells.Item(641,1).Value = 'Public, 4-year or above'
$ws.Cells.Item(641,2).Value = 'Public Property'
$ws.Cells.Item(641,3).Value = 'Criminal Offenses - Rape'
$ws.Cells.Item(641,4).Value = 'sum2015'
$ws.Cells.Item(641,5).Value = 71
$ws.Cells.Item(642,1).Value = 'Private nonprofit, 4-year or above'
$ws.Cells.Item(642,2).Value = 'Public Property'
$ws.Cells.Item(642,3).Value = 'Criminal Offenses - Rape'
$ws.Cells.Item(642,4).Value = 'sum2015'
$ws.Cells.Item(642,5).Value = 69
$ws.Cells.Item(643,1).Value = 'Private for-profit, 4-year or above'
$ws.Cells.Item(643,2).Value = 'Public Property'
$ws.Cells.Item(643,3).Value = 'Criminal Offenses - Rape'
$ws.Cells.Item(643,4).Value = 'sum2015'
$ws.Cells.Item(643,5).Value = 7
$ws.Cells.Item(644,1).Value = 'Public, 2-year'
$ws.Cells.Item(644,2).Value = 'Public Property'
$ws.Cells.Item(644,3).Value = 'Criminal Offenses - Rape'
$ws.Cells.Item(644,4).Value = 'sum2015'
$ws.Cells.Item(644,5).Value = 34
$ws.Cells.Item(645,1).Value = 'Private nonprofit, 2-year'
$ws.Cells.Item(645,2).Value = 'Public Property'
$ws.Cells.Item(645,3).Value = 'Criminal Offenses - Rape'
$ws.Cells.Item(645,4).Value = 'sum2015'
$ws.Cells.Item(645,5).Value = 0
$ws.Cells.Item(646,1).Value = 'Private for-profit, 2-year'
$ws.Cells.Item(646,2).Value = 'Public Property'
$ws.Cells.Item(646,3).Value = 'Criminal Offenses - Rape'
$ws.Cells.Item(646,4).Value = 'sum2015'
$ws.Cells.Item(646,5).Value = 8
$ws.Cells.Item(647,1).Value = 'Public, less-than 2-year'
$ws.Cells.Item(647,2).Value = 'Public Property'
$ws.Cells.Item(647,3).Value = 'Criminal Offenses - Rape'
$ws.Cells.Item(647,4).Value = 'sum2015'
$ws.Cells.Item(647,5).Value = 8
$ws.Cells.Item(648,1).Value = 'Private nonprofit, less-than 2-year'
$ws.Cells.Item(648,2).Value = 'Public Property'
$ws.Cells.Item(648,3).Value = 'Criminal Offenses - Rape'
$ws.Cells.Item(648,4).Value = 'sum2015'
$ws.Cells.Item(648,5).Value = 0
$ws.Cells.Item(649,1).Value = 'Private for-profit, less-than 2-year'
$ws.Cells.Item(649,2).Value = 'Public Property'
$ws.Cells.Item(649,3).Value = 'Criminal Offenses - Rape'
$ws.Cells.Item(649,4).Value = 'sum2015'
$ws.Cells.Item(649,5).Value = 8

# --- 2) Copy cell formatting from the last pre-existing row (577) onto the new block ---
#     so the new cells pick up the same style index (quotePrefix text style) as the
#     rest of the table.
$ws.Range("A577:E577").Copy()
$ws.Range("A578:E649").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- 3) Row 596 carries an extra formatted-but-empty cell in column F, matching the source data ---
$ws.Range("E577").Copy()
$ws.Range("F596").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- 4) Update the saved cursor/selection position to reflect where editing left off ---
[void]$ws.Range("D653").Select()
$aw = $excel.ActiveWindow
$aw.ScrollRow = 631
$aw.ScrollColumn = 1

# --- 5) Update the workbook window position recorded in the workbook view ---
try {
  $wb.Windows.Item(1).Left = 8520
} catch {}
